# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (values such as "1.00" or "0.200" would otherwise be coerced to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.762.93"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.338.08"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "0.668"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "236.82"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").Value = "72.37"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").Value = "57.07"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "31.77"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "7.13"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "2.686.86"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "16.23"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "0.882"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "2.336.05"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "43.628.25"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").Value = "0.0₃0997"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").Value = "76.30"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").Value = "254.10"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +22.31%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").Value = "10.53"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "22.42"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "174.29"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "0.134"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").Value = "0.0744"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("D36").Value = "5.10"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").Value = "3.68"
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "6.17"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").Value = "0.0273"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  +10.08%  "
$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.200"
$ws.Range("E43").Value = "  +6.92%  "
$ws.Range("D44").Value = "8.89"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "18.46"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("D46").Value = "59.54"
$ws.Range("E46").Value = "  +13.02%  "
$ws.Range("D47").Value = "4.66"
$ws.Range("E47").Value = "  +4.27%  "
$ws.Range("D48").Value = "2.44"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "1.22"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "98.83"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  -2.76%  "
